$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match details (columns F:V) between row pairs 14/15, 20/21 and
#    26/27. Columns A:E (index, country, tournament, season, match date)
#    stay attached to their original row.
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(14, 15),
    @(20, 21),
    @(26, 27)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("F${r1}:V${r1}")
    $range2 = $ws.Range("F${r2}:V${r2}")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# ---------------------------------------------------------------------------
# 2) Append three new match rows (43, 44, 45) at the bottom of the table.
#    Duplicate the formatting of the last existing row (42) first so the new
#    rows inherit the same styles/number formats, then overwrite the values.
# ---------------------------------------------------------------------------
$ws.Range("A42:V42").Copy($ws.Range("A43:V45"))

$newRows = @(
    @{
        Row = 43; A = 42; E = 45196.85416666666
        F = "Lugano"; G = 2; H = "Lausanne"; I = 1
        J = 1.81; K = "24/09/2023 22:12"; L = 2.2;  M = "27/09/2023 19:35"
        N = 3.89; O = "24/09/2023 22:12"; P = 3.84; Q = "27/09/2023 19:35"
        R = 3.93; S = "24/09/2023 22:12"; T = 3.13; U = "27/09/2023 19:35"
        V = "https://www.betexplorer.com/football/switzerland/super-league/lugano-lausanne/0tMuTbKB/"
    },
    @{
        Row = 44; A = 43; E = 45196.85416666666
        F = "Servette"; G = 2; H = "Winterthur"; I = 2
        J = 1.56; K = "24/09/2023 15:42"; L = 1.56; M = "27/09/2023 20:29"
        N = 4.47; O = "24/09/2023 15:42"; P = 4.44; Q = "27/09/2023 20:29"
        R = 4.91; S = "24/09/2023 15:42"; T = 5.75; U = "27/09/2023 20:29"
        V = "https://www.betexplorer.com/football/switzerland/super-league/servette-winterthur/vXKqSIZH/"
    },
    @{
        Row = 45; A = 44; E = 45196.85416666666
        F = "St. Gallen"; G = 2; H = "Young Boys"; I = 1
        J = 2.33; K = "24/09/2023 15:42"; L = 2.38; M = "27/09/2023 20:29"
        N = 3.72; O = "24/09/2023 15:42"; P = 3.68; Q = "27/09/2023 20:25"
        R = 2.75; S = "24/09/2023 15:42"; T = 2.93; U = "27/09/2023 20:29"
        V = "https://www.betexplorer.com/football/switzerland/super-league/st-gallen-young-boys/MeAlRxkO/"
    }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = "switzerland"
    $ws.Range("C$r").Value = "super-league"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
    $ws.Range("U$r").Value = $row.U
    $ws.Range("V$r").Value = $row.V
}
